$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "305.66"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.56%"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.45%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.041"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.34%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07965"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.47%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.929"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.95%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.776"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.32%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9197"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.06%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1314"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-5.39%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1903"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.75%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09037"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.62%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03420"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.14%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09843"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001412"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.04%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.006052"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "4.83%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.689"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.94%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.137"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.76%"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "12.58%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3445"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.13%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1342"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.09%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.169"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "5.20%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04422"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.30%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001233"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.92%"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-4.21%"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.59%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004448"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.15%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01930"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-4.14%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05328"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "8.25%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007599"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.54%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01014"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.58%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1354"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.73%"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "1.69%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009626"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-16.72%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006166"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-4.65%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.04%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "65.22"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.85%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001661"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "39.15%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.04%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.04%"
